$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B/C text changes (coin name / link swaps) ---
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"

# --- Column D price changes: values that are unambiguous text already (multiple dots) ---
$ws.Range("D2").Value = "42.472.27"
$ws.Range("D3").Value = "2.185.67"
$ws.Range("D14").Value = "2.510.30"
$ws.Range("D16").Value = "2.176.25"
$ws.Range("D18").Value = "42.403.23"

# --- Column D price changes: values that look numeric, force text format to avoid numeric coercion ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.88"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.586"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0917"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.101"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.773"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "170.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0807"
$ws.Range("D33").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0338"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.197"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "59.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.465"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0971"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.12"
$ws.Range("D50").Style = "Normal"

# --- Column E volume(1h) percentage changes ---
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("E5").Value = "  +5.49%  "
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("E7").Value = "  -2.17%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  -2.84%  "
$ws.Range("E10").Value = "  -2.30%  "
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("E14").Value = "  -1.54%  "
$ws.Range("E15").Value = "  -3.32%  "
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("E17").Value = "  -3.63%  "
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("E19").Value = "  -2.94%  "
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("E23").Value = "  -6.75%  "
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("E26").Value = "  -4.21%  "
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("E28").Value = "  -2.53%  "
$ws.Range("E29").Value = "  +7.41%  "
$ws.Range("E30").Value = "  -2.18%  "
$ws.Range("E31").Value = "  -1.18%  "
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("E33").Value = "  +1.87%  "
$ws.Range("E34").Value = "  -4.22%  "
$ws.Range("E35").Value = "  -1.18%  "
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("E37").Value = "  -3.33%  "
$ws.Range("E38").Value = "  +4.84%  "
$ws.Range("E39").Value = "  -6.22%  "
$ws.Range("E40").Value = "  -3.55%  "
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("E42").Value = "  -2.44%  "
$ws.Range("E43").Value = "  -6.49%  "
$ws.Range("E44").Value = "  +2.67%  "
$ws.Range("E45").Value = "  +6.74%  "
$ws.Range("E46").Value = "  +10.01%  "
$ws.Range("E47").Value = "  -1.35%  "
$ws.Range("E48").Value = "  -3.64%  "
$ws.Range("E49").Value = "  -1.22%  "
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("E51").Value = "  +0.49%  "
